$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.299.34"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.059.57"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.43"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  +2.39%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.84"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.07"
$ws.Range("E10").Value = "  +0.10%  "

$ws.Range("E11").Value = "  -0.26%  "

$ws.Range("E12").Value = "  +0.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.364.47"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.58"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.68"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.778"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("E17").Value = "  -1.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.060.80"
$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.218.05"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.34"
$ws.Range("E20").Value = "  +7.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.33"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0809"
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.20"
$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("E28").Value = "  +6.42%  "

$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.03"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0618"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.58"
$ws.Range("E35").Value = "  +5.49%  "

$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.67"
$ws.Range("E40").Value = "  -4.38%  "

$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.471.38"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0940"
$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.12"
$ws.Range("E44").Value = "  +1.87%  "

$ws.Range("E45").Value = "  +2.47%  "

$ws.Range("E46").Value = "  +3.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.19"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.12"
$ws.Range("E49").Value = "  -5.46%  "

$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("E51").Value = "  +0.82%  "
